$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string reorder: swap "Montserrat" / "Islas Malvinas" ---
# Before: A213 = Montserrat, A214 = Islas Malvinas
# After:  A213 = Islas Malvinas, A214 = Montserrat
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 15:38"

# --- Updated per-country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 5613357
$ws.Range("C4").Value = 1330
$ws.Range("D4").Value = 2974788
$ws.Range("E4").Value = 2464791
$ws.Range("G4").Value = 62
$ws.Range("H4").Value = 173778
$ws.Range("B6").Value = 2724642
$ws.Range("C6").Value = 23038
$ws.Range("D6").Value = 1997918
$ws.Range("E6").Value = 674571
$ws.Range("G6").Value = 228
$ws.Range("H6").Value = 52153
$ws.Range("B16").Value = 301323
$ws.Range("C16").Value = 1409
$ws.Range("D16").Value = 272911
$ws.Range("E16").Value = 24942
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 3470
$ws.Range("B36").Value = 85219
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 5790
$ws.Range("B45").Value = 63973
$ws.Range("C45").Value = 489
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 6175
$ws.Range("B61").Value = 36352
$ws.Range("C61").Value = 650
$ws.Range("D61").Value = 32062
$ws.Range("E61").Value = 4048
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 242
$ws.Range("B67").Value = 29890
$ws.Range("C67").Value = 108
$ws.Range("D67").Value = 27597
$ws.Range("E67").Value = 1612
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 681
$ws.Range("D72").Value = 14929
$ws.Range("E72").Value = 8406
$ws.Range("B83").Value = 12970
$ws.Range("C83").Value = 130
$ws.Range("D83").Value = 9513
$ws.Range("E83").Value = 2908
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 549
$ws.Range("E87").Value = 941
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 262
$ws.Range("B96").Value = 8131
$ws.Range("C96").Value = 32
$ws.Range("D96").Value = 6935
$ws.Range("E96").Value = 1131
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 65
$ws.Range("D158").Value = 526
$ws.Range("E158").Value = 437
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 26
$ws.Range("B194").Value = 97
$ws.Range("C194").Value = 3
$ws.Range("D194").Value = 88
$ws.Range("E194").Value = 8
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
